{"js": "const replacements = [\n  { old: \"379\u00f76=63, 1\", new: \"838\u00f78=104, 6\" },\n  { old: \"845\u00f76=140, 5\", new: \"244\u00f76=40, 4\" },\n  { old: \"165\u00f76=27, 3\", new: \"504\u00f78=63, 0\" },\n  { old: \"640\u00f78=80, 0\", new: \"695\u00f77=99, 2\" },\n  { old: \"991\u00f77=141, 4\", new: \"507\u00f79=56, 3\" },\n  { old: \"143\u00f77=20, 3\", new: \"131\u00f75=26, 1\" },\n  { old: \"751\u00f73=250, 1\", new: \"183\u00f78=22, 7\" },\n  { old: \"474\u00f73=158, 0\", new: \"337\u00f73=112, 1\" },\n  { old: \"284\u00f74=71, 0\", new: \"828\u00f79=92, 0\" },\n  { old: \"809\u00f79=89, 8\", new: \"230\u00f78=28, 6\" },\n  { old: \"369\u00f74=92, 1\", new: \"978\u00f73=326, 0\" },\n  { old: \"823\u00f77=117, 4\", new: \"395\u00f76=65, 5\" },\n  { old: \"797\u00f72=398, 1\", new: \"441\u00f75=88, 1\" },\n  { old: \"646\u00f75=129, 1\", new: \"381\u00f75=76, 1\" },\n  { old: \"118\u00f78=14, 6\", new: \"301\u00f79=33, 4\" },\n  { old: \"710\u00f78=88, 6\", new: \"168\u00f77=24, 0\" },\n  { old: \"223\u00f72=111, 1\", new: \"180\u00f73=60, 0\" },\n  { old: \"225\u00f72=112, 1\", new: \"514\u00f73=171, 1\" },\n  { old: \"572\u00f72=286, 0\", new: \"506\u00f72=253, 0\" },\n  { old: \"103\u00f73=34, 1\", new: \"510\u00f74=127, 2\" },\n  { old: \"870\u00f75=174, 0\", new: \"574\u00f76=95, 4\" },\n  { old: \"214\u00f78=26, 6\", new: \"408\u00f76=68, 0\" },\n  { old: \"463\u00f76=77, 1\", new: \"230\u00f73=76, 2\" },\n  { old: \"130\u00f74=32, 2\", new: \"504\u00f79=56, 0\" },\n  { old: \"260\u00f78=32, 4\", new: \"513\u00f78=64, 1\" },\n];\n\nconst body = context.document.body;\n\nfor (const { old, new: newText } of replacements) {\n  const results = body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"379\u00f76=63, 1\"; new=\"838\u00f78=104, 6\"},\n    @{old=\"845\u00f76=140, 5\"; new=\"244\u00f76=40, 4\"},\n    @{old=\"165\u00f76=27, 3\"; new=\"504\u00f78=63, 0\"},\n    @{old=\"640\u00f78=80, 0\"; new=\"695\u00f77=99, 2\"},\n    @{old=\"991\u00f77=141, 4\"; new=\"507\u00f79=56, 3\"},\n    @{old=\"143\u00f77=20, 3\"; new=\"131\u00f75=26, 1\"},\n    @{old=\"751\u00f73=250, 1\"; new=\"183\u00f78=22, 7\"},\n    @{old=\"474\u00f73=158, 0\"; new=\"337\u00f73=112, 1\"},\n    @{old=\"284\u00f74=71, 0\"; new=\"828\u00f79=92, 0\"},\n    @{old=\"809\u00f79=89, 8\"; new=\"230\u00f78=28, 6\"},\n    @{old=\"369\u00f74=92, 1\"; new=\"978\u00f73=326, 0\"},\n    @{old=\"823\u00f77=117, 4\"; new=\"395\u00f76=65, 5\"},\n    @{old=\"797\u00f72=398, 1\"; new=\"441\u00f75=88, 1\"},\n    @{old=\"646\u00f75=129, 1\"; new=\"381\u00f75=76, 1\"},\n    @{old=\"118\u00f78=14, 6\"; new=\"301\u00f79=33, 4\"},\n    @{old=\"710\u00f78=88, 6\"; new=\"168\u00f77=24, 0\"},\n    @{old=\"223\u00f72=111, 1\"; new=\"180\u00f73=60, 0\"},\n    @{old=\"225\u00f72=112, 1\"; new=\"514\u00f73=171, 1\"},\n    @{old=\"572\u00f72=286, 0\"; new=\"506\u00f72=253, 0\"},\n    @{old=\"103\u00f73=34, 1\"; new=\"510\u00f74=127, 2\"},\n    @{old=\"870\u00f75=174, 0\"; new=\"574\u00f76=95, 4\"},\n    @{old=\"214\u00f78=26, 6\"; new=\"408\u00f76=68, 0\"},\n    @{old=\"463\u00f76=77, 1\"; new=\"230\u00f73=76, 2\"},\n    @{old=\"130\u00f74=32, 2\"; new=\"504\u00f79=56, 0\"},\n    @{old=\"260\u00f78=32, 4\"; new=\"513\u00f78=64, 1\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
